$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "7667 ms"
$ws.Range("E3").Value = "5561 ms"
$ws.Range("E4").Value = "6727 ms"
$ws.Range("E5").Value = "9108 ms"
$ws.Range("E6").Value = "6003 ms"
$ws.Range("E7").Value = "7754 ms"
$ws.Range("E8").Value = "7132 ms"
$ws.Range("E9").Value = "6596 ms"
